# Add the "ReadOnly property of constant" check as a new row (row 33) on
# the "Workflow" sheet, mirroring the formatting of the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workflow")

# Duplicate the formatting (styles, borders, alignment) of the last existing
# data row (32) down onto the new row (33) so the new row renders like its
# neighbours.
$ws.Range("A32:G32").Copy()
$ws.Range("A33:G33").PasteSpecial(-4122)
$ws.Rows.Item(33).RowHeight = 87

# Column A (Run) - reuse existing "No" shared string.
$ws.Range("A33").Value = "No"

# Column C (Check Filename) then B (Issue) - new shared strings, inserted in
# this order so they land at shared-string indices 168 and 169 respectively.
$ws.Range("C33").Value = "Checks\Custom\ReadOnlyPropertyOfConstant.xaml"
$ws.Range("B33").Value = "ReadOnly property of constant"

# Column D (Arguments) - new shared string (index 170).
$ws.Range("D33").Value = "{ " + [char]10 + "`"NamingPattern`" : `"(^([A-Z][A-Z0-9]*)`$)`"" + [char]10 + "}"

# Column E (Action) - reuse existing "Fix" shared string.
$ws.Range("E33").Value = "Fix"

# Column F (Explanation) - new shared string (index 171).
$ws.Range("F33").Value = "Constant is a variable that does not change the value set by the default value. " + [char]10 + "It must follow a specific naming convention so that it can be distinguished from other variables. In addition, ReadOnly must be specified in the variable's Modifiers property."

# Column G (Suggestion) - new shared string (index 172).
$ws.Range("G33").Value = "Variables that do not have their values modified should follow the naming convention of constants and have ReadOnly checked in their Modifiers property."

# Extend the data-validation ranges so row 33 is covered the same way row 32
# was: the "Yes, No" list on column A, and the "Fix, Double check" list on
# column E (E2:E6 and E11:E32).
$ws.Range("E11:E32").Validation.Delete()
$ws.Range("A2:A33").Validation.Delete()
$ws.Range("A2:A33").Validation.Add(3, 1, 1, """Yes, No""")
$ws.Range("E11:E33").Validation.Add(3, 1, 1, """Fix, Double check""")
